$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing B2 value
$ws.Range("B2").Value = 1299

# Update existing A3/B3 row, then add new rows 4 and 5
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 1203

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 1157

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 589

# Copy style (format) from A3 (which already had the bold/border style) to the new A4 and A5 cells
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
